# Auto-generated Excel COM-interop script to update Siren_Profits market price data
# Applies scheduled-runner refreshed values for columns H-N across ALC, ARM, BSM, CRP, CUL, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1017.2963
$ws.Range("I19").Value = 462.6154
$ws.Range("K19").Value = 462.6154
$ws.Range("M19").Value = -287.6154
$ws.Range("H80").Value = 250672.75
$ws.Range("L80").Value = 2466.75
$ws.Range("M80").Value = -1500571.75
$ws.Range("I80").Value = 500523.25
$ws.Range("N80").Value = -4462.75
$ws.Range("J80").Value = 822.25
$ws.Range("K80").Value = 1501569.75
$ws.Range("I83").Value = 500523.25
$ws.Range("H83").Value = 250672.75
$ws.Range("M83").Value = -4499717.25
$ws.Range("J83").Value = 822.25
$ws.Range("K83").Value = 4504709.25
$ws.Range("L83").Value = 7400.25
$ws.Range("N83").Value = -17384.25
$ws.Range("K86").Value = 3775.6667
$ws.Range("I86").Value = 3775.6667
$ws.Range("H86").Value = 71462480
$ws.Range("M86").Value = -2652.6667
$ws.Range("I88").Value = 0
$ws.Range("M88").Value = ""
$ws.Range("N88").Value = -1612
$ws.Range("H88").Value = 800
$ws.Range("L88").Value = 800
$ws.Range("K88").Value = 0
$ws.Range("J88").Value = 800
$ws.Range("I89").Value = 3775.6667
$ws.Range("H89").Value = 71462480
$ws.Range("K89").Value = 18878.3335
$ws.Range("M89").Value = -13262.3335
$ws.Range("H91").Value = 800
$ws.Range("L91").Value = 800
$ws.Range("I91").Value = 0
$ws.Range("N91").Value = -3608
$ws.Range("K91").Value = 0
$ws.Range("M91").Value = ""
$ws.Range("J91").Value = 800
$ws.Range("I94").Value = 500000000
$ws.Range("K94").Value = 500000000
$ws.Range("H94").Value = 100252000
$ws.Range("M94").Value = -499999549
$ws.Range("H129").Value = 1320.8
$ws.Range("I129").Value = 1022.25
$ws.Range("K129").Value = 3066.75
$ws.Range("M129").Value = 1933.25
$ws.Range("I138").Value = 1367.4445
$ws.Range("N138").Value = -23194.2938
$ws.Range("L138").Value = 12914.2938
$ws.Range("J138").Value = 4304.7646
$ws.Range("H138").Value = 3689.9768
$ws.Range("K138").Value = 4102.333500000001
$ws.Range("M138").Value = 1037.666499999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I45").Value = 4115.2
$ws.Range("L45").Value = 6148.625
$ws.Range("J45").Value = 6148.625
$ws.Range("H45").Value = 4822.478
$ws.Range("M45").Value = -3738.2
$ws.Range("N45").Value = -6902.625
$ws.Range("K45").Value = 4115.2
$ws.Range("I74").Value = 1375.8298
$ws.Range("K74").Value = 1375.8298
$ws.Range("H74").Value = 2502.8103
$ws.Range("M74").Value = -501.8298
$ws.Range("K77").Value = 6879.148999999999
$ws.Range("I77").Value = 1375.8298
$ws.Range("M77").Value = -2511.148999999999
$ws.Range("H77").Value = 2502.8103
$ws.Range("J122").Value = 2594841
$ws.Range("K122").Value = 9649.5
$ws.Range("M122").Value = -7199.5
$ws.Range("H122").Value = 1262005.6
$ws.Range("I122").Value = 3216.5
$ws.Range("N122").Value = -7789423
$ws.Range("L122").Value = 7784523
$ws.Range("H132").Value = 1512.2222
$ws.Range("I132").Value = 912.1951
$ws.Range("M132").Value = -206.5853000000002
$ws.Range("K132").Value = 2736.5853

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 89718.45
$ws.Range("L82").Value = 103100.336
$ws.Range("J82").Value = 103100.336
$ws.Range("N82").Value = -103866.336
$ws.Range("L85").Value = 103100.336
$ws.Range("N85").Value = -105752.336
$ws.Range("H85").Value = 89718.45
$ws.Range("J85").Value = 103100.336
$ws.Range("M99").Value = -11144.454
$ws.Range("K99").Value = 12642.454
$ws.Range("I99").Value = 12642.454
$ws.Range("H99").Value = 11476.15
$ws.Range("I105").Value = 1513.75
$ws.Range("M105").Value = 233.25
$ws.Range("H105").Value = 3191.7273
$ws.Range("K105").Value = 1513.75
$ws.Range("H107").Value = 4164.8096
$ws.Range("K107").Value = 3905.1667
$ws.Range("L107").Value = 4511
$ws.Range("N107").Value = -8351
$ws.Range("J107").Value = 4511
$ws.Range("M107").Value = -1985.1667
$ws.Range("I107").Value = 3905.1667

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I94").Value = 2209.818
$ws.Range("K94").Value = 2209.818
$ws.Range("H94").Value = 2192.3333
$ws.Range("M94").Value = -1758.818
$ws.Range("I105").Value = 11764.777
$ws.Range("M105").Value = -10017.777
$ws.Range("H105").Value = 8914.076999999999
$ws.Range("K105").Value = 11764.777
$ws.Range("H132").Value = 14448.371
$ws.Range("I132").Value = 2957.484
$ws.Range("M132").Value = -6342.451999999999
$ws.Range("K132").Value = 8872.451999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M98").Value = -1854.8
$ws.Range("K98").Value = 3352.8
$ws.Range("H98").Value = 1367.6923
$ws.Range("I98").Value = 1117.6
$ws.Range("M99").Value = -3436
$ws.Range("K99").Value = 5682
$ws.Range("I99").Value = 1894
$ws.Range("H99").Value = 6310
$ws.Range("M103").Value = -77908.5
$ws.Range("L103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("H103").Value = 26262.5
$ws.Range("K103").Value = 78787.5
$ws.Range("I103").Value = 26262.5
$ws.Range("N103").Value = ""
$ws.Range("H104").Value = 6499.5
$ws.Range("M104").Value = -6376
$ws.Range("K104").Value = 8997
$ws.Range("I104").Value = 2999
$ws.Range("H107").Value = 532.4375
$ws.Range("L107").Value = 1915.7499
$ws.Range("N107").Value = -5755.7499
$ws.Range("J107").Value = 638.5833
$ws.Range("I108").Value = 83333630
$ws.Range("K108").Value = 250000890
$ws.Range("H108").Value = 55558388
$ws.Range("M108").Value = -249998010
$ws.Range("H109").Value = 2185.4167
$ws.Range("M109").Value = -5725
$ws.Range("I109").Value = 2255
$ws.Range("K109").Value = 6765
$ws.Range("N110").Value = -37820
$ws.Range("J110").Value = 9880
$ws.Range("H110").Value = 23953.6
$ws.Range("L110").Value = 29640
$ws.Range("K111").Value = 11088.6
$ws.Range("H111").Value = 3696.2
$ws.Range("M111").Value = -8021.599999999999
$ws.Range("L111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("N111").Value = ""
$ws.Range("I111").Value = 3696.2
$ws.Range("I112").Value = 3158.2856
$ws.Range("J112").Value = 7560
$ws.Range("K112").Value = 9474.856800000001
$ws.Range("H112").Value = 4992.3335
$ws.Range("N112").Value = -24896
$ws.Range("M112").Value = -8366.856800000001
$ws.Range("L112").Value = 22680
$ws.Range("L113").Value = 13699.0005
$ws.Range("N113").Value = -18039.0005
$ws.Range("H113").Value = 4566.3335
$ws.Range("J113").Value = 4566.3335
$ws.Range("I114").Value = 849.5
$ws.Range("L114").Value = 11296.5
$ws.Range("H114").Value = 2307.5
$ws.Range("J114").Value = 3765.5
$ws.Range("N114").Value = -17804.5
$ws.Range("K114").Value = 2548.5
$ws.Range("M114").Value = 705.5
$ws.Range("H116").Value = 125002690
$ws.Range("I116").Value = 142859940
$ws.Range("K116").Value = 428579820
$ws.Range("M116").Value = -428576378
$ws.Range("H117").Value = 2019.875
$ws.Range("L117").Value = 2400
$ws.Range("J117").Value = 800
$ws.Range("N117").Value = -9284

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3049.8572
$ws.Range("M82").Value = -5479.5
$ws.Range("L82").Value = 1933.6
$ws.Range("J82").Value = 1933.6
$ws.Range("N82").Value = -2655.6
$ws.Range("K82").Value = 5840.5
$ws.Range("I82").Value = 5840.5
$ws.Range("L85").Value = 1933.6
$ws.Range("M85").Value = -4592.5
$ws.Range("K85").Value = 5840.5
$ws.Range("N85").Value = -4429.6
$ws.Range("H85").Value = 3049.8572
$ws.Range("J85").Value = 1933.6
$ws.Range("I85").Value = 5840.5
$ws.Range("K136").Value = 6708.999899999999
$ws.Range("L136").Value = 28843.386
$ws.Range("M136").Value = -4158.999899999999
$ws.Range("J136").Value = 9614.462
$ws.Range("N136").Value = -33943.386
$ws.Range("H136").Value = 7284.5264
$ws.Range("I136").Value = 2236.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11057.5
$ws.Range("L41").Value = 9000.75
$ws.Range("N41").Value = -9780.75
$ws.Range("J41").Value = 9000.75
$ws.Range("M62").Value = -254176.27
$ws.Range("I62").Value = 254800.27
$ws.Range("H62").Value = 160499.92
$ws.Range("K62").Value = 254800.27
$ws.Range("K65").Value = 1274001.35
$ws.Range("H65").Value = 160499.92
$ws.Range("M65").Value = -1270881.35
$ws.Range("I65").Value = 254800.27
$ws.Range("K136").Value = 6288.75
$ws.Range("L136").Value = 14998.9995
$ws.Range("M136").Value = -3738.75
$ws.Range("J136").Value = 4999.6665
$ws.Range("N136").Value = -20098.9995
$ws.Range("H136").Value = 2888.0908
$ws.Range("I136").Value = 2096.25
